$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Combinations" labels (column B), replacing the old ones
$combinations = @(
    "(3, 'distance', 'chebyshev', 'brute', 70)",
    "(3, 'distance', 'chebyshev', 'brute', 150)",
    "(3, 'distance', 'chebyshev', 'brute', 30)",
    "(3, 'distance', 'chebyshev', 'brute', 300)",
    "(3, 'distance', 'chebyshev', 'brute', 10)",
    "(3, 'distance', 'manhattan', 'kd_tree', 10)",
    "(3, 'distance', 'manhattan', 'brute', 10)",
    "(3, 'distance', 'manhattan', 'kd_tree', 30)",
    "(3, 'distance', 'manhattan', 'ball_tree', 30)",
    "(3, 'distance', 'manhattan', 'ball_tree', 10)"
)

# New metric values: R2 (col C), MSE (col D), MAE (col E)
# Rows 2-6 share one set of values, rows 7-11 share another
$r2Top    = 0.9996211712091984
$mseTop   = 31001.32435897436
$maeTop   = 25.11153846153846

$r2Bottom  = 0.9995973586037076
$mseBottom = 32408.03076923077
$maeBottom = 21.28974358974359

for ($i = 0; $i -lt $combinations.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $combinations[$i]

    if ($row -le 6) {
        $r2 = $r2Top
        $mse = $mseTop
        $mae = $maeTop
    } else {
        $r2 = $r2Bottom
        $mse = $mseBottom
        $mae = $maeBottom
    }

    $ws.Cells.Item($row, 3).Value = $r2
    $ws.Cells.Item($row, 4).Value = $mse
    $ws.Cells.Item($row, 5).Value = $mae
}
